$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.855.23"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.413.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.67"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.89"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.204"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +11.38%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.72"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "688.60"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.72"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.965.24"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.867.07"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.418.42"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.37"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.03%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.80"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.55%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.01"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.66"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +8.98%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "558.07"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.84%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.77"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.12%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.672.61"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.12"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0742"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +9.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.31"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.72"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +6.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.342"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.37"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.59%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.62"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.70"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.90%  "
